# Bump the published "term" value set to 1.1.0 and refresh its publication
# date, matching the FHIR metadata table on the "Metadata" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B3").Value = "1.1.0"
$ws.Range("B8").Value = "2023-07-10T23:08:03+02:00"

# The header/body cell styles already carried a top-aligned, wrap-text
# <alignment> element, but it was never actually switched on because the
# xf records were missing applyAlignment="true". Turn wrapping + top
# vertical alignment on (its already-intended values) for every sheet so
# the alignment formatting that was defined becomes effective, across
# both worksheets in the workbook.
foreach ($sheet in $wb.Worksheets) {
    $used = $sheet.UsedRange
    $used.WrapText = $true
    $used.VerticalAlignment = -4160
}
